$d = $word.ActiveDocument

# Locate the paragraph that contains "Test avec ... -set-upstream." using a
# wildcard Find so we don't depend on a hard-coded paragraph index.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Test avec*set-upstream.", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target paragraph 'Test avec ... -set-upstream.'"
}

$targetStart = $rng.Start
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -le $targetStart -and $p.Range.End -gt $targetStart) {
        $target = $p
        break
    }
}

if ($null -eq $target) {
    throw "Could not resolve paragraph object for found range"
}

# Replace that single paragraph with two paragraphs:
#  1. The original sentence, now reading "Test avec git push -set-upstream."
#     split across several runs/proofErr tags (as Word would do after typing
#     "git push" in the middle of the sentence, with proofing language
#     switched to English for the inserted text).
#  2. A brand-new paragraph "Ok ca marche je test le git push simple." that
#     inherits the bullet/paragraph formatting and carries forward the
#     _GoBack bookmark (tracking the last edit location).
$xml = '<w:p w14:paraId="47E6EC2F" w14:textId="3322428D" w:rsidR="0001454E" w:rsidRDefault="0001454E" w:rsidP="005039A2"><w:pPr><w:pStyle w:val="paragraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="45"/></w:numPr><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:jc w:val="both"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji"/><w:lang w:val="en-US"/></w:rPr><w:t>Test avec</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji"/><w:lang w:val="en-US"/></w:rPr><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> push</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> &#8211;set-upstream.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="paragraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="45"/></w:numPr><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:jc w:val="both"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji"/></w:rPr><w:t>Ok &#231;a marche je test le git push simple.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$null = $target.Range.InsertXML($xml)
